$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2, 3 and 4 get cyclically rotated for columns
# A, B, E, F, G, H, Q, R:
#   new row2 = old row3
#   new row3 = old row4
#   new row4 = old row2
# (columns C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
# stay untouched because they are identical across the three rows.)

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values before overwriting anything.
$orig2 = @{}
$orig3 = @{}
$orig4 = @{}
foreach ($col in $cols) {
    $orig2[$col] = $ws.Range("$col`2").Value2
    $orig3[$col] = $ws.Range("$col`3").Value2
    $orig4[$col] = $ws.Range("$col`4").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col`2").Value2 = $orig3[$col]
    $ws.Range("$col`3").Value2 = $orig4[$col]
    $ws.Range("$col`4").Value2 = $orig2[$col]
}
